$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.27285266666667
$ws.Range("H2").Value = 63.818558
$ws.Range("I2").Value = 0.1403158092128609
$ws.Range("J2").Value = 0.1403158092128609
$ws.Range("M2").Value = 80.019516
$ws.Range("N2").Value = 240.058548
$ws.Range("O2").Value = 0.3109587206407793
$ws.Range("P2").Value = 0.3109587206407793
$ws.Range("Q2").Value = 1702.243374325976
$ws.Range("R2").Value = 15320.19036893378
$ws.Range("S2").Value = 0.04363242451850691
$ws.Range("T2").Value = 0.04363242451850691
$ws.Range("G3").Value = 21.27285266666667
$ws.Range("H3").Value = 63.818558
$ws.Range("I3").Value = 0.1403158092128609
$ws.Range("J3").Value = 0.1403158092128609
$ws.Range("O3").Value = 0.1376395003539564
$ws.Range("P3").Value = 0.1376395003539564
$ws.Range("Q3").Value = 753.4631189640105
$ws.Range("R3").Value = 6781.168070676094
$ws.Range("S3").Value = 0.01931299787181925
$ws.Range("T3").Value = 0.01931299787181925
$ws.Range("G4").Value = 21.27285266666667
$ws.Range("H4").Value = 63.818558
$ws.Range("I4").Value = 0.1403158092128609
$ws.Range("J4").Value = 0.1403158092128609
$ws.Range("M4").Value = 14.174735
$ws.Range("N4").Value = 42.524205
$ws.Range("O4").Value = 0.05508353063547745
$ws.Range("P4").Value = 0.05508353063547745
$ws.Range("Q4").Value = 301.5370492440433
$ws.Range("R4").Value = 2713.83344319639
$ws.Range("S4").Value = 0.007729090175418434
$ws.Range("T4").Value = 0.007729090175418434
$ws.Range("G5").Value = 21.27285266666667
$ws.Range("H5").Value = 63.818558
$ws.Range("I5").Value = 0.1403158092128609
$ws.Range("J5").Value = 0.1403158092128609
$ws.Range("M5").Value = 127.7183863333333
$ws.Range("N5").Value = 383.155159
$ws.Range("O5").Value = 0.4963182483697869
$ws.Range("P5").Value = 0.4963182483697869
$ws.Range("Q5").Value = 2716.934415293414
$ws.Range("R5").Value = 24452.40973764072
$ws.Range("S5").Value = 0.06964129664711635
$ws.Range("T5").Value = 0.06964129664711635
$ws.Range("I6").Value = 0.06281023381453604
$ws.Range("J6").Value = 0.06281023381453604
$ws.Range("M6").Value = 80.019516
$ws.Range("N6").Value = 240.058548
$ws.Range("O6").Value = 0.3109587206407793
$ws.Range("P6").Value = 0.3109587206407793
$ws.Range("Q6").Value = 761.9833071586601
$ws.Range("R6").Value = 6857.84976442794
$ws.Range("S6").Value = 0.01953138995011635
$ws.Range("T6").Value = 0.01953138995011635
$ws.Range("I7").Value = 0.06281023381453604
$ws.Range("J7").Value = 0.06281023381453604
$ws.Range("O7").Value = 0.1376395003539564
$ws.Range("P7").Value = 0.1376395003539564
$ws.Range("S7").Value = 0.008645169199347916
$ws.Range("T7").Value = 0.008645169199347916
$ws.Range("I8").Value = 0.06281023381453604
$ws.Range("J8").Value = 0.06281023381453604
$ws.Range("M8").Value = 14.174735
$ws.Range("N8").Value = 42.524205
$ws.Range("O8").Value = 0.05508353063547745
$ws.Range("P8").Value = 0.05508353063547745
$ws.Range("Q8").Value = 134.9784651708917
$ws.Range("R8").Value = 1214.806186538025
$ws.Range("S8").Value = 0.003459809438544497
$ws.Range("T8").Value = 0.003459809438544497
$ws.Range("I9").Value = 0.06281023381453604
$ws.Range("J9").Value = 0.06281023381453604
$ws.Range("M9").Value = 127.7183863333333
$ws.Range("N9").Value = 383.155159
$ws.Range("O9").Value = 0.4963182483697869
$ws.Range("P9").Value = 0.4963182483697869
$ws.Range("Q9").Value = 1216.1942894436
$ws.Range("R9").Value = 10945.7486049924
$ws.Range("S9").Value = 0.03117386522652729
$ws.Range("T9").Value = 0.03117386522652729
$ws.Range("G10").Value = 1.192664333333333
$ws.Range("H10").Value = 3.577993
$ws.Range("I10").Value = 0.007866818036737088
$ws.Range("J10").Value = 0.007866818036737088
$ws.Range("M10").Value = 80.019516
$ws.Range("N10").Value = 240.058548
$ws.Range("O10").Value = 0.3109587206407793
$ws.Range("P10").Value = 0.3109587206407793
$ws.Range("Q10").Value = 95.436422703796
$ws.Range("R10").Value = 858.927804334164
$ws.Range("S10").Value = 0.002446255672217572
$ws.Range("T10").Value = 0.002446255672217572
$ws.Range("G11").Value = 1.192664333333333
$ws.Range("H11").Value = 3.577993
$ws.Range("I11").Value = 0.007866818036737088
$ws.Range("J11").Value = 0.007866818036737088
$ws.Range("O11").Value = 0.1376395003539564
$ws.Range("P11").Value = 0.1376395003539564
$ws.Range("Q11").Value = 42.24297523944989
$ws.Range("R11").Value = 380.186777155049
$ws.Range("S11").Value = 0.001082784903951985
$ws.Range("T11").Value = 0.001082784903951985
$ws.Range("G12").Value = 1.192664333333333
$ws.Range("H12").Value = 3.577993
$ws.Range("I12").Value = 0.007866818036737088
$ws.Range("J12").Value = 0.007866818036737088
$ws.Range("M12").Value = 14.174735
$ws.Range("N12").Value = 42.524205
$ws.Range("O12").Value = 0.05508353063547745
$ws.Range("P12").Value = 0.05508353063547745
$ws.Range("Q12").Value = 16.90570086895167
$ws.Range("R12").Value = 152.151307820565
$ws.Range("S12").Value = 0.000433332112330334
$ws.Range("T12").Value = 0.000433332112330334
$ws.Range("G13").Value = 1.192664333333333
$ws.Range("H13").Value = 3.577993
$ws.Range("I13").Value = 0.007866818036737088
$ws.Range("J13").Value = 0.007866818036737088
$ws.Range("M13").Value = 127.7183863333333
$ws.Range("N13").Value = 383.155159
$ws.Range("O13").Value = 0.4963182483697869
$ws.Range("P13").Value = 0.4963182483697869
$ws.Range("Q13").Value = 152.3251640906541
$ws.Range("R13").Value = 1370.926476815887
$ws.Range("S13").Value = 0.003904445348237198
$ws.Range("T13").Value = 0.003904445348237198
$ws.Range("G14").Value = 119.6189703333333
$ws.Range("H14").Value = 358.856911
$ws.Range("I14").Value = 0.7890071389358658
$ws.Range("J14").Value = 0.7890071389358659
$ws.Range("M14").Value = 80.019516
$ws.Range("N14").Value = 240.058548
$ws.Range("O14").Value = 0.3109587206407793
$ws.Range("P14").Value = 0.3109587206407793
$ws.Range("Q14").Value = 9571.852110491691
$ws.Range("R14").Value = 86146.66899442521
$ws.Range("S14").Value = 0.2453486504999385
$ws.Range("T14").Value = 0.2453486504999385
$ws.Range("G15").Value = 119.6189703333333
$ws.Range("H15").Value = 358.856911
$ws.Range("I15").Value = 0.7890071389358658
$ws.Range("J15").Value = 0.7890071389358659
$ws.Range("O15").Value = 0.1376395003539564
$ws.Range("P15").Value = 0.1376395003539564
$ws.Range("Q15").Value = 4236.784031125402
$ws.Range("R15").Value = 38131.05628012862
$ws.Range("S15").Value = 0.1085985483788372
$ws.Range("T15").Value = 0.1085985483788372
$ws.Range("G16").Value = 119.6189703333333
$ws.Range("H16").Value = 358.856911
$ws.Range("I16").Value = 0.7890071389358658
$ws.Range("J16").Value = 0.7890071389358659
$ws.Range("M16").Value = 14.174735
$ws.Range("N16").Value = 42.524205
$ws.Range("O16").Value = 0.05508353063547745
$ws.Range("P16").Value = 0.05508353063547745
$ws.Range("Q16").Value = 1695.567205447862
$ws.Range("R16").Value = 15260.10484903075
$ws.Range("S16").Value = 0.04346129890918417
$ws.Range("T16").Value = 0.04346129890918418
$ws.Range("G17").Value = 119.6189703333333
$ws.Range("H17").Value = 358.856911
$ws.Range("I17").Value = 0.7890071389358658
$ws.Range("J17").Value = 0.7890071389358659
$ws.Range("M17").Value = 127.7183863333333
$ws.Range("N17").Value = 383.155159
$ws.Range("O17").Value = 0.4963182483697869
$ws.Range("P17").Value = 0.4963182483697869
$ws.Range("Q17").Value = 15277.54186582821
$ws.Range("R17").Value = 137497.8767924538
$ws.Range("S17").Value = 0.391598641147906
$ws.Range("T17").Value = 0.391598641147906
